# Insert a new data row at row 18 (shifts existing rows 18..110 down to 19..111,
# matching every subsequent row's content staying identical but moved down by one).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(18).Insert()

# Populate the newly inserted row 18 with its record.
$ws.Range("A18").Value = 3
$ws.Range("B18").Value = "Femacal de La Calera"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = 44462
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100101
$ws.Range("H18").Value = "Berries"
$ws.Range("I18").Value = 100112025
$ws.Range("J18").Value = "Frutilla"
$ws.Range("K18").Value = "Sin especificar"
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 17000
$ws.Range("O18").Value = 17000
$ws.Range("P18").Value = 17000
$ws.Range("Q18").Value = '$/bandeja 7 kilos'
$ws.Range("R18").Value = "Provincia de Melipilla"
$ws.Range("S18").Value = 2429
$ws.Range("T18").Value = 7
